$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 344, shifting existing rows 344:357 down to 345:358
$ws.Rows(344).Insert()

# Populate the newly inserted row 344 with the new weekly record
$ws.Cells.Item(344, 1).Value = 8
$ws.Cells.Item(344, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(344, 3).Value = "Coquimbo"
$ws.Cells.Item(344, 4).Value = 44939
$ws.Cells.Item(344, 5).Value = 4
$ws.Cells.Item(344, 6).Value = 100112021
$ws.Cells.Item(344, 7).Value = "Ají"
$ws.Cells.Item(344, 8).Value = "Inferno"
$ws.Cells.Item(344, 9).Value = "Primera"
$ws.Cells.Item(344, 10).Value = 440
$ws.Cells.Item(344, 11).Value = 14000
$ws.Cells.Item(344, 12).Value = 15000
$ws.Cells.Item(344, 13).Value = 14500
$ws.Cells.Item(344, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(344, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(344, 16).Value = 967
$ws.Cells.Item(344, 17).Value = 15
$ws.Cells.Item(344, 18).Value = "Hortaliza"
